# repull data, push all data, mean calculation
# Update column F ("dSF") values for specific rows to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 0
    6  = -1
    9  = -6
    10 = -4
    12 = -2
    13 = 1
    15 = -2
    19 = 1
    28 = -1
    37 = 1
    40 = -3
    41 = 4
    42 = -2
    44 = 1
    51 = 2
    59 = 0
    61 = -2
    70 = -3
    75 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
